# "both pages and sidebar"
#
# 1. Normalize the Verdict column (col B) wording:
#      PASSED         -> Passed
#      FAILED         -> Failed
#      Not_applicable -> Not_Applicable
# 2. Move the view/selection: scroll back to the top of the sheet and
#    select D36 (instead of the old scrolled-down B20 selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Text -eq "PASSED") {
        $cell.Value = "Passed"
    } elseif ($cell.Text -eq "FAILED") {
        $cell.Value = "Failed"
    } elseif ($cell.Text -eq "Not_applicable") {
        $cell.Value = "Not_Applicable"
    }
}

# Scroll the window back to the top-left (A1) ...
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# ... and move the active selection to D36.
$ws.Range("D36").Select()
